# Insert a new data row at row 443 (shifting the existing rows 443-508 down
# to 444-509), and populate the new row with its own data, exactly as
# described by the target diff:
#   - dimension A1:R508 -> A1:R509
#   - new row 443 inserted with fresh values
#   - all previously existing rows 443..508 shift down to 444..509 unchanged

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 443..508 down to 444..509 by inserting a blank row at 443.
$ws.Rows.Item(443).Insert()

# Populate the newly inserted row 443 with the new record's data.
$ws.Cells.Item(443, 1).Value = 5
$ws.Cells.Item(443, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(443, 3).Value = "Maule"
$ws.Cells.Item(443, 4).Value = 44984
$ws.Cells.Item(443, 5).Value = 7
$ws.Cells.Item(443, 6).Value = 100112032
$ws.Cells.Item(443, 7).Value = "Zapallo italiano"
$ws.Cells.Item(443, 8).Value = "Sin especificar"
$ws.Cells.Item(443, 9).Value = "Primera"
$ws.Cells.Item(443, 10).Value = 400
$ws.Cells.Item(443, 11).Value = 4000
$ws.Cells.Item(443, 12).Value = 4000
$ws.Cells.Item(443, 13).Value = 4000
$ws.Cells.Item(443, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(443, 15).Value = "Región del Maule"
$ws.Cells.Item(443, 16).Value = 80
$ws.Cells.Item(443, 17).Value = 50
$ws.Cells.Item(443, 18).Value = "Hortaliza"
